## Replace the " m:self.name " Word FIELD (fldChar begin/instrText.../fldChar end)
## with plain literal text runs: "{" "m" ":" "self"(colored) ".name}"
## i.e. turn the field code into literal M2Doc template braces text.

$d = $word.ActiveDocument

# Locate the field whose code contains "m:self.name" (robust to paragraph
# position -- does not assume a hard-coded paragraph index).
$target = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields($i)
    if ($candidate.Code.Text -match "self") {
        $target = $candidate
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'm:self.name' field."
}

# The whole paragraph that holds the field (begin fldChar .. end fldChar).
$fieldParagraphRange = $target.Code.Paragraphs(1).Range

# Rebuild that paragraph's content as literal text runs (no more field codes):
#   {  m  :  self(colored)  .name}
$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>{</w:t>
            </w:r>
            <w:r>
              <w:t>m</w:t>
            </w:r>
            <w:r>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
              </w:rPr>
              <w:t>self</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">.name}</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$fieldParagraphRange.InsertXML($newContentXml)

Write-Output "Replaced field with literal template text."
